$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grades")

# Replace "Maths" with "Math" across column B (Subject column)
$ws.Columns.Item(2).Replace("Maths", "Math", 1, 1, $false, $false, $false)

# Update the selection to match the post-edit state (column B selected)
$ws.Columns.Item(2).Select()
